$wb = $excel.ActiveWorkbook

# Rename Sheet1 -> Fueltrade
$wsFuel = $wb.Worksheets.Item("Sheet1")
$wsFuel.Name = "Fueltrade"

$wsFuel.Range("B2").Value = "~TradeLinks_DINS"
$wsFuel.Range("B3").Value = "Reg1"
$wsFuel.Range("C3").Value = "Reg2"
$wsFuel.Range("D3").Value = "Comm"
$wsFuel.Range("E3").Value = "Comm1"
$wsFuel.Range("F3").Value = "Comm2"
$wsFuel.Range("G3").Value = "Tech"
$wsFuel.Range("H3").Value = "TradeLink"

$data = @(
  @("DKISLBH","DKE","H2","H2","H2","TB_H2_DKISLBH_DKE_01","U"),
  @("DKISLBH","DKE","H2","H2","H2","TB_H2_DKISLBH_DKE_02","U"),
  @("DKISL1","DKW","H2","H2","H2","TB_H2_DKISL1_DKW_01","U"),
  @("DKISL1","DKW","H2","H2","H2","TB_H2_DKISL1_DKW_02","U"),
  @("DKISL2","DKW","H2","H2","H2","TB_H2_DKISL2_DKW_01","U"),
  @("DKISL2","DKW","H2","H2","H2","TB_H2_DKISL2_DKW_02","U"),
  @("DKISL3","DKW","H2","H2","H2","TB_H2_DKISL3_DKW_01","U"),
  @("DKISL3","DKW","H2","H2","H2","TB_H2_DKISL3_DKW_02","U")
)

$r = 4
foreach ($row in $data) {
  $wsFuel.Cells.Item($r, 2).Value = $row[0]
  $wsFuel.Cells.Item($r, 3).Value = $row[1]
  $wsFuel.Cells.Item($r, 4).Value = $row[2]
  $wsFuel.Cells.Item($r, 5).Value = $row[3]
  $wsFuel.Cells.Item($r, 6).Value = $row[4]
  $wsFuel.Cells.Item($r, 7).Value = $row[5]
  $wsFuel.Cells.Item($r, 8).Value = $row[6]
  $r++
}

$wsFuel.Range("G12").Select()

# BI sheet: clear C12, move selection
$wsBI = $wb.Worksheets.Item("BI")
$wsBI.Range("C12").ClearContents()
$wsBI.Range("I20").Select()

# Uni sheet: add selection
$wsUni = $wb.Worksheets.Item("Uni")
$wsUni.Range("C35").Select()
